# Requirements Progress.xlsx - "Error handling for parsing an incomplete decision tree"
#
# Summary of the change (from the "Developmnet PB" sheet):
#   - Row 26 ("Error handling for tree that can't be parsed") is marked Done
#     and highlighted green.
#   - "include equipment in tree" is re-cased to "Include equipment in tree",
#     marked Done and highlighted green (moves conceptually next to the
#     error-handling item).
#   - "Add loading things for ajax in tree editor" is reworded to
#     "Add loading gif for ajax in tree editor" and its status drops back to
#     Not Started.
#   - A handful of other in-flight items (28, 31, 36, 41) are marked Done,
#     a batch of "In Progress" items (29,30,32-35,37,39,40,42-45) revert to
#     "Not Started", and "Complete Report" moves from N/A to In Progress.
#   - A brand-new backlog item "Find out what browsers are compatible" is
#     appended as row 46 (Not Started).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Developmnet PB")

$GREEN = 5287936  # BGR encoding of RGB(0,176,80) - same green used elsewhere on the sheet

# --- Row 26: status -> Done, highlight green -------------------------------
$ws.Range("C26").Value = "Done"
$ws.Range("A26").Interior.Color = $GREEN

# --- Row 27: reword + status -> Done, highlight green ----------------------
$ws.Range("A27").Value = "Include equipment in tree"
$ws.Range("C27").Value = "Done"
$ws.Range("A27").Interior.Color = $GREEN

# --- Row 28: status -> Done -------------------------------------------------
$ws.Range("C28").Value = "Done"

# --- Row 29: reword + status -> Not Started ---------------------------------
$ws.Range("A29").Value = "Add loading gif for ajax in tree editor"
$ws.Range("C29").Value = "Not Started"

# --- Row 30: status -> Not Started ------------------------------------------
$ws.Range("C30").Value = "Not Started"

# --- Row 31: status -> Done -------------------------------------------------
$ws.Range("C31").Value = "Done"

# --- Rows 32-35: status -> Not Started --------------------------------------
$ws.Range("C32").Value = "Not Started"
$ws.Range("C33").Value = "Not Started"
$ws.Range("C34").Value = "Not Started"
$ws.Range("C35").Value = "Not Started"

# --- Row 36: status -> Done -------------------------------------------------
$ws.Range("C36").Value = "Done"

# --- Row 37: status -> Not Started ------------------------------------------
$ws.Range("C37").Value = "Not Started"

# --- Row 38: status -> In Progress ------------------------------------------
$ws.Range("C38").Value = "In Progress"

# --- Rows 39-40: status -> Not Started --------------------------------------
$ws.Range("C39").Value = "Not Started"
$ws.Range("C40").Value = "Not Started"

# --- Row 41: status -> Done -------------------------------------------------
$ws.Range("C41").Value = "Done"

# --- Rows 42-45: status -> Not Started --------------------------------------
$ws.Range("C42").Value = "Not Started"
$ws.Range("C43").Value = "Not Started"
$ws.Range("C44").Value = "Not Started"
$ws.Range("C45").Value = "Not Started"

# --- Row 46 (new): append backlog item --------------------------------------
$ws.Range("A46").Value = "Find out what browsers are compatible"
$ws.Range("A46").WrapText = $true
$ws.Range("B46").Interior.Color = 49407  # same orange used for the other "Not Started" rows in column B
$ws.Range("C46").Value = "Not Started"

# --- Best-effort view-state update (scroll position / active selection) ----
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 18
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("D40").Select()
